$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "convocatoria_ocds_id"
$ws.Range("B1").Value = "item_etapa_id"
$ws.Range("C1").Value = "item_id"
$ws.Range("D1").Value = "item_descripcion"
$ws.Range("E1").Value = "item_clasificacion"
$ws.Range("F1").Value = "item_cantidad"
$ws.Range("G1").Value = "item_unidad"

# New header cells F1/G1 should carry the same bold/border/center-top style as
# the rest of the header row, so copy the formatting from an existing header cell.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data row (row 2) ---
$ws.Range("A2").Value = "ocds-twb234-0005"

# item_etapa_id (B2) has no value in the new dataset row, so make sure the
# cell stays empty/absent rather than keeping a stale blank entry.
$ws.Range("B2").ClearContents()

# item_id (C2) must be stored as text "3245", not converted to a number.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "3245"
$ws.Range("C2").ClearFormats()

$ws.Range("D2").Value = "Servicio de consultoria"
$ws.Range("E2").Value = "Soporte y mantenimiento de hardware"

# item_cantidad (F2) must be stored as text "1", not converted to a number.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"
$ws.Range("F2").ClearFormats()

$ws.Range("G2").Value = "Cantidad"
